# Update the "duration (minutes)" column (H) for the easy aerobic run (WT = "E")
# rows of microcycles 2-10 in the Training Plan sheet, bumping the easy run
# duration up in line with the plan's new progression.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Plan")

$updates = @{
    25 = 44; 26 = 44; 28 = 44; 29 = 44; 31 = 44
    35 = 48; 36 = 48; 38 = 48; 39 = 48; 41 = 48
    45 = 51; 46 = 51; 48 = 51; 49 = 51; 51 = 51
    57 = 55; 60 = 55; 61 = 55; 62 = 55
    65 = 59; 66 = 59; 68 = 59; 69 = 59; 71 = 59
    75 = 62; 76 = 62; 78 = 62; 79 = 62; 81 = 62
    85 = 66; 86 = 66; 88 = 66; 89 = 66; 91 = 66
    97 = 70; 100 = 70; 101 = 70; 102 = 70
    105 = 50; 106 = 50; 108 = 50; 109 = 50; 111 = 50
    117 = 30; 120 = 30; 121 = 30; 122 = 30
}

foreach ($row in $updates.Keys) {
    $ws.Range("H$row").Value = $updates[$row]
}
